# Apply the authored edit to the presentation:
#   1. Bump the cached "datetimeFigureOut" date field text from
#      13/09/2022 to 20/09/2022 on the Slide Master and every Slide Layout
#      (the placeholder that normally auto-fills the footer/date area).
#   2. On slide 4, change the red "14" label (TextBox 57) to "12".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    foreach ($sh in $shapes) {
        $isDatePlaceholder = $false
        try {
            if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "13/09/2022") {
                $sh.TextFrame.TextRange.Text = "20/09/2022"
            }
        }
    }
}

# Slide Master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout's date placeholder.
foreach ($layout in $master.CustomLayouts) {
    Update-DatePlaceholder $layout.Shapes
}

# Slide 4: the red "14" callout becomes "12".
$slide4 = $p.Slides.Item(4)
foreach ($sh in $slide4.Shapes) {
    if ($sh.Name -eq "TextBox 57") {
        if ($sh.TextFrame.TextRange.Text -eq "14") {
            $sh.TextFrame.TextRange.Text = "12"
        }
    }
}
